$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update row 2 values and delete row 3 ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("A2").Value = 46040.14583333334
$schedule.Range("B2").Value = 46040.8125
$schedule.Range("C2").Value = 16
$schedule.Range("D2").Value = 60.48
$schedule.Range("E2").Value = 295.28607225
$schedule.Range("F2").Value = 4.882375533234127
$schedule.Rows.Item(3).Delete()

# --- Sheet "Detailed": update individual cells per diff ---
$detailed = $wb.Worksheets.Item("Detailed")
$detailed.Range("B9").Value = 57.06022
$detailed.Range("E9").Value = "ON"
$detailed.Range("E10").Value = "ON"
$detailed.Range("B11").Value = 36.06
$detailed.Range("C11").Value = "historical"
$detailed.Range("E11").Value = "ON"
$detailed.Range("B12").Value = 36.2
$detailed.Range("C12").Value = "historical"
$detailed.Range("E12").Value = "ON"
$detailed.Range("B13").Value = 56.98
$detailed.Range("E13").Value = "ON"
$detailed.Range("E14").Value = "ON"
$detailed.Range("B15").Value = 36.06
$detailed.Range("E15").Value = "ON"
$detailed.Range("B16").Value = 8.058479999999999
$detailed.Range("B19").Value = 0.009549999999999999
$detailed.Range("B20").Value = 0
$detailed.Range("B21").Value = -0.9066
$detailed.Range("B22").Value = -4.72132
$detailed.Range("B23").Value = 0.66745
$detailed.Range("B24").Value = 0.7
$detailed.Range("B25").Value = -0.93546
$detailed.Range("B26").Value = -4.77113
$detailed.Range("B27").Value = -4.61593
$detailed.Range("B28").Value = -5.51
$detailed.Range("B29").Value = -6.4985
$detailed.Range("B30").Value = -10
$detailed.Range("B31").Value = -20.71215
$detailed.Range("B32").Value = -11.45471
$detailed.Range("B33").Value = -11.01
$detailed.Range("B34").Value = -7.05469
$detailed.Range("B35").Value = -6.56987
$detailed.Range("B37").Value = 0.66377
$detailed.Range("B38").Value = 3.98566
$detailed.Range("B39").Value = 16.85794
$detailed.Range("B40").Value = 41.08477
$detailed.Range("E40").Value = "ON"
$detailed.Range("B41").Value = 57.08646
$detailed.Range("E42").Value = "OFF"
$detailed.Range("B43").Value = 56.98
$detailed.Range("E43").Value = "OFF"
$detailed.Range("E44").Value = "OFF"
$detailed.Range("B45").Value = 46.33506
$detailed.Range("E45").Value = "OFF"
$detailed.Range("B46").Value = 56.33763
$detailed.Range("E46").Value = "OFF"
$detailed.Range("E47").Value = "OFF"
$detailed.Range("E48").Value = "OFF"
$detailed.Range("E49").Value = "OFF"
